# Resize/reposition the right-hand "Notes" textbox + its body text box on the
# two "Two-column slide" slides (19 and 20) so both columns share the same
# left edge and width.

$p = $ppt.ActivePresentation

$newLeft  = 4648200 / 12700       # 366.0 pt
# 4495799 EMU isn't an exact multiple of 12700 (it is 1 EMU short of 354.0pt);
# Shape.Width is a single-precision COM property, so feed it a point value
# that still rounds back to the exact target EMU count on save.
$newWidth = 353.99994

# --- Slide 19: "Two-column slide (60% / 40%)" ---
$s19 = $p.Slides.Item(19)

$notes19 = $s19.Shapes.Item("TextBox 5")
$notes19.Left  = $newLeft
$notes19.Width = $newWidth

$body19 = $s19.Shapes.Item("TextBox 6")
$body19.Left  = $newLeft
$body19.Width = $newWidth

# --- Slide 20: "Two-column slide (Auto + default)" ---
$s20 = $p.Slides.Item(20)

$notes20 = $s20.Shapes.Item("TextBox 5")
$notes20.Left  = $newLeft
$notes20.Width = $newWidth

$body20 = $s20.Shapes.Item("TextBox 6")
$body20.Left   = $newLeft
$body20.Width  = $newWidth
$body20.Height = 342900 / 12700   # 27.0 pt (grows from 21.6 pt)
